# "case distr insert ok" - add a new "distr/case/insert" panel-description
# block (label/textbox/richtextbox/button controls) as rows 140-143 on
# Plan1 (sheet1), mirroring the existing repeated blocks further up the
# sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows 128-134 (the "D=28" group) already carry the alternate-section
# fill/border/alignment formatting that the new "D=30" group should reuse.
# Stamp that existing formatting onto the four new rows first so no new
# style entries get minted - the new rows end up sharing the very same
# style as that block, just like in the source file.
$ws.Range("A128:E128").Copy()
$ws.Range("A140:E143").PasteSpecial(-4122)

# Now fill in the actual cell values. The order below controls the order
# new entries are appended to the shared-string table, matching how the
# original edit built up the rows.
$ws.Range("A140").Value = "label112;"
$ws.Range("C140").Value = "textBox44;"
$ws.Range("C141").Value = "richTextBox7;"
$ws.Range("A141").Value = "label113;"
$ws.Range("E140").Value = "distr/case/insert"
$ws.Range("A142").Value = "button42;"
$ws.Range("B140").Value = "título"
$ws.Range("A143").Value = "button43;"

$ws.Range("B141").Value = "Description"
$ws.Range("B142").Value = "Submit"
$ws.Range("C142").Value = "-"
$ws.Range("B143").Value = "clear"
$ws.Range("C143").Value = "-"

$ws.Range("D140").Value = 30
$ws.Range("D141").Value = 30
$ws.Range("D142").Value = 30
$ws.Range("D143").Value = 30

$ws.Range("E141").Value = "distr/case/insert"
$ws.Range("E142").Value = "distr/case/insert"
$ws.Range("E143").Value = "distr/case/insert"

# Scroll/select to match the author's final view position.
$win = $excel.ActiveWindow
$win.ScrollRow = 131
$win.ScrollColumn = 1
$ws.Range("C140").Select()
